$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 6494
$wsExpo.Range("F15").Value = 3193
$wsExpo.Range("F18").Value = 1864

# Sheet "全部类型" (all types) - same events, different row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6494
$wsAll.Range("F16").Value = 3193
$wsAll.Range("F19").Value = 1864
